$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "62.202.23"
Set-TextValue $ws.Range("E2") "  +0.74%  "

Set-TextValue $ws.Range("D3") "2.415.93"
Set-TextValue $ws.Range("E3") "  +1.42%  "

Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.14%  "

Set-TextValue $ws.Range("D5") "561.34"
Set-TextValue $ws.Range("E5") "  +1.64%  "

Set-TextValue $ws.Range("D6") "142.99"

Set-TextValue $ws.Range("E7") "  +0.04%  "

Set-TextValue $ws.Range("E8") "  +1.23%  "

Set-TextValue $ws.Range("D9") "2.412.09"
Set-TextValue $ws.Range("E9") "  +1.11%  "

Set-TextValue $ws.Range("E10") "  +0.32%  "

Set-TextValue $ws.Range("E11") "  -2.06%  "

Set-TextValue $ws.Range("E12") "  -0.94%  "

Set-TextValue $ws.Range("E13") "  -0.49%  "

Set-TextValue $ws.Range("D14") "25.66"
Set-TextValue $ws.Range("E14") "  -0.53%  "

Set-TextValue $ws.Range("E15") "  +0.04%  "

Set-TextValue $ws.Range("D16") "2.859.86"
Set-TextValue $ws.Range("E16") "  +1.68%  "

Set-TextValue $ws.Range("D17") "62.018.16"
Set-TextValue $ws.Range("E17") "  +0.67%  "

Set-TextValue $ws.Range("D18") "2.412.62"
Set-TextValue $ws.Range("E18") "  +1.19%  "

Set-TextValue $ws.Range("D19") "11.28"
Set-TextValue $ws.Range("E19") "  +2.40%  "

Set-TextValue $ws.Range("E20") "  -0.04%  "

Set-TextValue $ws.Range("D21") "323.29"
Set-TextValue $ws.Range("E21") "  -0.01%  "

Set-TextValue $ws.Range("E22") "  +1.76%  "

Set-TextValue $ws.Range("D24") "65.71"
Set-TextValue $ws.Range("E24") "  +2.02%  "

Set-TextValue $ws.Range("E25") "  -3.86%  "

Set-TextValue $ws.Range("D26") "9.01"
Set-TextValue $ws.Range("E26") "  -0.02%  "

Set-TextValue $ws.Range("D27") "579.06"
Set-TextValue $ws.Range("E27") "  +6.60%  "

Set-TextValue $ws.Range("B28") "Binance-PegBSC-USD"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D28") "1.00"
Set-TextValue $ws.Range("E28") "  +0.57%  "

Set-TextValue $ws.Range("B29") "WrappedeETH"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D29") "2.535.42"
Set-TextValue $ws.Range("E29") "  +1.40%  "

Set-TextValue $ws.Range("D30") "0.0₃0944"
Set-TextValue $ws.Range("E30") "  +2.00%  "

Set-TextValue $ws.Range("D31") "8.22"
Set-TextValue $ws.Range("E31") "  -0.98%  "

Set-TextValue $ws.Range("E32") "  +1.51%  "

Set-TextValue $ws.Range("E33") "  +0.18%  "

Set-TextValue $ws.Range("E34") "  +0.79%  "

Set-TextValue $ws.Range("E35") "  -0.28%  "

Set-TextValue $ws.Range("E36") "  -0.03%  "

Set-TextValue $ws.Range("E37") "  -1.74%  "

Set-TextValue $ws.Range("D38") "4.74"
Set-TextValue $ws.Range("E38") "  -0.66%  "

Set-TextValue $ws.Range("E39") "  +0.50%  "

Set-TextValue $ws.Range("D40") "152.48"
Set-TextValue $ws.Range("E40") "  +3.97%  "

Set-TextValue $ws.Range("D41") "18.64"
Set-TextValue $ws.Range("E41") "  +0.26%  "

Set-TextValue $ws.Range("E42") "  -6.34%  "

Set-TextValue $ws.Range("D43") "0.997"
Set-TextValue $ws.Range("E43") "  -0.30%  "

Set-TextValue $ws.Range("E44") "  +1.50%  "

Set-TextValue $ws.Range("D45") "149.01"
Set-TextValue $ws.Range("E45") "  +0.73%  "

Set-TextValue $ws.Range("E46") "  +0.84%  "

Set-TextValue $ws.Range("D47") "0.0534"
Set-TextValue $ws.Range("E47") "  +0.75%  "

Set-TextValue $ws.Range("D48") "20.08"
Set-TextValue $ws.Range("E48") "  -0.76%  "

Set-TextValue $ws.Range("D49") "0.594"
Set-TextValue $ws.Range("E49") "  +1.39%  "

Set-TextValue $ws.Range("E50") "  +1.07%  "

Set-TextValue $ws.Range("E51") "  +1.20%  "
